# Add "Age" and "Phone_Number" as new optional columns to the static
# data-upload template (Sample_Person_Add.xlsx), per commit:
# "Age and phone number included in static data upload. optional"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cells: H1 = Age, I1 = Phone_Number
$ws.Cells.Item(1, 8).Value = "Age"
$ws.Cells.Item(1, 9).Value = "Phone_Number"

# Give the new Phone_Number column a sensible custom width, matching the
# other header columns' explicit custom widths.
$ws.Columns.Item(9).ColumnWidth = 14.85

# Reflect the author's final cursor position in the sheet view.
$ws.Range("G11").Select()
